$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.40130033333333
$ws.Range("H2").Value = 109.203901
$ws.Range("I2").Value = 0.1897437225523226
$ws.Range("J2").Value = 0.1897437225523226
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 73.78963405897177
$ws.Range("R2").Value = 664.106706530746
$ws.Range("S2").Value = 0.001251603588456957
$ws.Range("T2").Value = 0.001251603588456957
$ws.Range("G3").Value = 36.40130033333333
$ws.Range("H3").Value = 109.203901
$ws.Range("I3").Value = 0.1897437225523226
$ws.Range("J3").Value = 0.1897437225523226
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 9334.9061772117
$ws.Range("R3").Value = 84014.15559490531
$ws.Range("S3").Value = 0.1583366311312751
$ws.Range("T3").Value = 0.1583366311312751
$ws.Range("G4").Value = 36.40130033333333
$ws.Range("H4").Value = 109.203901
$ws.Range("I4").Value = 0.1897437225523226
$ws.Range("J4").Value = 0.1897437225523226
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 1777.849178892122
$ws.Range("R4").Value = 16000.6426100291
$ws.Range("S4").Value = 0.03015548783259059
$ws.Range("T4").Value = 0.03015548783259058
$ws.Range("I5").Value = 0.6107553255746098
$ws.Range("J5").Value = 0.6107553255746098
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 237.5172752357649
$ws.Range("R5").Value = 2137.655477121884
$ws.Range("S5").Value = 0.00402871592733502
$ws.Range("T5").Value = 0.004028715927335019
$ws.Range("I6").Value = 0.6107553255746098
$ws.Range("J6").Value = 0.6107553255746098
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("Q6").Value = 30047.60096819119
$ws.Range("R6").Value = 270428.4087137207
$ws.Range("S6").Value = 0.5096608172125537
$ws.Range("T6").Value = 0.5096608172125536
$ws.Range("I7").Value = 0.6107553255746098
$ws.Range("J7").Value = 0.6107553255746098
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 5722.618063305828
$ws.Range("R7").Value = 51503.56256975244
$ws.Range("S7").Value = 0.09706579243472109
$ws.Range("T7").Value = 0.09706579243472105
$ws.Range("G8").Value = 38.27317166666666
$ws.Range("H8").Value = 114.819515
$ws.Range("I8").Value = 0.1995009518730676
$ws.Range("J8").Value = 0.1995009518730676
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 77.58413314079887
$ws.Range("R8").Value = 698.25719826719
$ws.Range("S8").Value = 0.001315965049626638
$ws.Range("T8").Value = 0.001315965049626638
$ws.Range("G9").Value = 38.27317166666666
$ws.Range("H9").Value = 114.819515
$ws.Range("I9").Value = 0.1995009518730676
$ws.Range("J9").Value = 0.1995009518730676
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 9814.936920961747
$ws.Range("R9").Value = 88334.43228865573
$ws.Range("S9").Value = 0.1664788073205087
$ws.Range("T9").Value = 0.1664788073205087
$ws.Range("G10").Value = 38.27317166666666
$ws.Range("H10").Value = 114.819515
$ws.Range("I10").Value = 0.1995009518730676
$ws.Range("J10").Value = 0.1995009518730676
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 1869.27187210594
$ws.Range("R10").Value = 16823.44684895346
$ws.Range("S10").Value = 0.03170617950293234
$ws.Range("T10").Value = 0.03170617950293234
